# The post "「連休初日」" (previously row 169) was removed from the sheet.
# Deleting the entire row 169 shifts every row below it up by one,
# which matches the rest of the diff (rows 170-372 becoming 169-371)
# and shrinks the used range from A1:C372 to A1:C371.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("169:169").Delete()
